# Fruta / hortaliza, semanal
# Insert a new weekly record at row 120 (shifting the existing rows 120-198
# down to 121-199) and populate the new row with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 120; this pushes the existing
# data (previously rows 120-198) down to rows 121-199, preserving it
# unchanged.
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with the new weekly record.
$ws.Cells.Item(120, 1).Value = 4
$ws.Cells.Item(120, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(120, 3).Value = "Los Lagos"
$ws.Cells.Item(120, 4).Value = 44574
$ws.Cells.Item(120, 5).Value = 10
$ws.Cells.Item(120, 6).Value = 100112021
$ws.Cells.Item(120, 7).Value = "Ají"
$ws.Cells.Item(120, 8).Value = "Inferno"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 70
$ws.Cells.Item(120, 11).Value = 21000
$ws.Cells.Item(120, 12).Value = 21000
$ws.Cells.Item(120, 13).Value = 21000
$ws.Cells.Item(120, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(120, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(120, 16).Value = 1750
$ws.Cells.Item(120, 17).Value = 12
$ws.Cells.Item(120, 18).Value = "Hortaliza"

# Apply the same date-number format used by the rest of column D (style
# index 2 in the original workbook) to the new row's date cell.
$ws.Cells.Item(120, 4).NumberFormat = $ws.Cells.Item(121, 4).NumberFormat
